# This sheet is a weekly NFL picks worksheet. Each of rows 2-5 holds one
# game: column A is the fixed game label, B is the favored team (user
# input), C is the point spread (user input) and D is the over/under
# (user input). Columns E-J and the helper columns (L:X) recompute
# automatically from B/C/D (including two What-If "Data Table" result
# ranges, G2:G5 and J2:J5, driven by the row-input cells N14/U14).
#
# The edit: swap in this week's games/lines -
#   - Game 1 (row 2): favorite NE -> LAC, spread -4 -> -6, o/u 34 -> 40.5
#   - Game 2 (row 3): spread -9 -> -6.5, o/u 44 -> 42 (favorite stays KC)
#   - Game 3 (row 4): favorite/spread/o-u cleared out (no game entered)
#   - Game 4 (row 5): favorite LAR -> CLE, spread -1 -> -4, o/u 46 -> 39.5
#
# Note: row 5 (CLE) is entered before row 2 (LAC) so the new shared
# strings land in the same order ("CLE" then "LAC") as the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Game 4 (row 5): NE's LAR opponent becomes CLE
$ws.Range("B5").Value2 = "CLE"
$ws.Range("C5").Value2 = -4
$ws.Range("D5").Value2 = 39.5

# Game 1 (row 2): NE becomes LAC
$ws.Range("B2").Value2 = "LAC"
$ws.Range("C2").Value2 = -6
$ws.Range("D2").Value2 = 40.5

# Game 2 (row 3): KC line updated, favorite unchanged
$ws.Range("C3").Value2 = -6.5
$ws.Range("D3").Value2 = 42

# Game 3 (row 4): no game this slot anymore - clear favorite/spread/o-u
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
